# Auto-generated edit script
# Updates market price / profit columns (H-N) across all 8 job sheets
# per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 389.09756
$ws.Range("J17").Value = 389.09756
$ws.Range("L17").Value = 1167.29268
$ws.Range("N17").Value = -1503.29268
$ws.Range("H21").Value = 15337.333
$ws.Range("I21").Value = 15337.333
$ws.Range("K21").Value = 15337.333
$ws.Range("M21").Value = -14869.333
$ws.Range("H23").Value = 15337.333
$ws.Range("I23").Value = 15337.333
$ws.Range("K23").Value = 15337.333
$ws.Range("M23").Value = -15103.333
$ws.Range("H34").Value = 4647.5713
$ws.Range("I34").Value = 4647.5713
$ws.Range("K34").Value = 4647.5713
$ws.Range("M34").Value = -4444.5713
$ws.Range("H36").Value = 4647.5713
$ws.Range("I36").Value = 4647.5713
$ws.Range("K36").Value = 4647.5713
$ws.Range("M36").Value = -3932.5713
$ws.Range("H86").Value = 4091.75
$ws.Range("I86").Value = 3937.2
$ws.Range("K86").Value = 3937.2
$ws.Range("M86").Value = -2814.2
$ws.Range("H89").Value = 4091.75
$ws.Range("I89").Value = 3937.2
$ws.Range("K89").Value = 19686
$ws.Range("M89").Value = -14070
$ws.Range("H96").Value = 2941.3333
$ws.Range("I96").Value = 1287.5
$ws.Range("J96").Value = 6249
$ws.Range("K96").Value = 3862.5
$ws.Range("L96").Value = 18747
$ws.Range("M96").Value = -2489.5
$ws.Range("N96").Value = -21493
$ws.Range("H125").Value = 2298.8333
$ws.Range("I125").Value = 2974.5
$ws.Range("J125").Value = 1961
$ws.Range("K125").Value = 26770.5
$ws.Range("L125").Value = 17649
$ws.Range("M125").Value = -24310.5
$ws.Range("N125").Value = -22569
$ws.Range("H139").Value = 94447.5
$ws.Range("J139").Value = 94447.5
$ws.Range("L139").Value = 94447.5
$ws.Range("N139").Value = -104727.5
$ws.Range("H140").Value = 119948.5
$ws.Range("J140").Value = 119948.5
$ws.Range("L140").Value = 119948.5
$ws.Range("N140").Value = -130308.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 53937.477
$ws.Range("I32").Value = 62876.35
$ws.Range("K32").Value = 62876.35
$ws.Range("M32").Value = -62589.35
$ws.Range("H61").Value = 2086.3547
$ws.Range("I61").Value = 2006.0741
$ws.Range("K61").Value = 2006.0741
$ws.Range("M61").Value = -1794.0741
$ws.Range("H132").Value = 33154.938
$ws.Range("I132").Value = 39002.184
$ws.Range("K132").Value = 117006.552
$ws.Range("M132").Value = -114476.552
$ws.Range("H136").Value = 2086.3547
$ws.Range("I136").Value = 2006.0741
$ws.Range("K136").Value = 6018.2223
$ws.Range("M136").Value = -3468.2223

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H105").Value = 3540.375
$ws.Range("I105").Value = 3348.5
$ws.Range("K105").Value = 3348.5
$ws.Range("M105").Value = -1601.5
$ws.Range("H134").Value = 2240.0344
$ws.Range("I134").Value = 2240.0344
$ws.Range("K134").Value = 6720.1032
$ws.Range("M134").Value = -4185.1032
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1055.4
$ws.Range("I105").Value = 1055.4
$ws.Range("K105").Value = 1055.4
$ws.Range("M105").Value = 691.5999999999999
$ws.Range("H107").Value = 2463.4666
$ws.Range("I107").Value = 483.54544
$ws.Range("J107").Value = 3609.7368
$ws.Range("K107").Value = 483.54544
$ws.Range("L107").Value = 3609.7368
$ws.Range("M107").Value = 1436.45456
$ws.Range("N107").Value = -7449.736800000001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 64378.3
$ws.Range("J37").Value = 64378.3
$ws.Range("L37").Value = 193134.9
$ws.Range("N37").Value = -193358.9
$ws.Range("H133").Value = 19127.8
$ws.Range("I133").Value = 17659.75
$ws.Range("J133").Value = 25000
$ws.Range("K133").Value = 52979.25
$ws.Range("L133").Value = 75000
$ws.Range("M133").Value = -47919.25
$ws.Range("N133").Value = -85120
$ws.Range("H139").Value = 1222.5
$ws.Range("I139").Value = 1222.5
$ws.Range("K139").Value = 3667.5
$ws.Range("M139").Value = 1472.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1809.6111
$ws.Range("I80").Value = 1038.2
$ws.Range("K80").Value = 1038.2
$ws.Range("M80").Value = -40.20000000000005
$ws.Range("H83").Value = 1809.6111
$ws.Range("I83").Value = 1038.2
$ws.Range("K83").Value = 5191
$ws.Range("M83").Value = -199
$ws.Range("H122").Value = 2697.5
$ws.Range("I122").Value = 1973.0358
$ws.Range("J122").Value = 6078.3335
$ws.Range("K122").Value = 5919.107400000001
$ws.Range("L122").Value = 18235.0005
$ws.Range("M122").Value = -3469.107400000001
$ws.Range("N122").Value = -23135.0005

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5059.8
$ws.Range("I7").Value = 3824.75
$ws.Range("K7").Value = 3824.75
$ws.Range("M7").Value = -3712.75
$ws.Range("H22").Value = 70505.75
$ws.Range("I22").Value = 139549.38
$ws.Range("J22").Value = 1462.125
$ws.Range("K22").Value = 139549.38
$ws.Range("L22").Value = 1462.125
$ws.Range("M22").Value = -139254.38
$ws.Range("N22").Value = -2052.125
$ws.Range("H27").Value = 70505.75
$ws.Range("I27").Value = 139549.38
$ws.Range("J27").Value = 1462.125
$ws.Range("K27").Value = 139549.38
$ws.Range("L27").Value = 1462.125
$ws.Range("M27").Value = -139442.38
$ws.Range("N27").Value = -1676.125
$ws.Range("H126").Value = 5059.8
$ws.Range("I126").Value = 3824.75
$ws.Range("K126").Value = 11474.25
$ws.Range("M126").Value = -9004.25

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 50000
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H113").Value = 1596.1538
$ws.Range("I113").Value = 574.875
$ws.Range("K113").Value = 1724.625
$ws.Range("M113").Value = 445.375
$ws.Range("H122").Value = 5341.44
$ws.Range("I122").Value = 5706.227
$ws.Range("K122").Value = 17118.681
$ws.Range("M122").Value = -14668.681
$ws.Range("H136").Value = 2751.3333
$ws.Range("J136").Value = 3400
$ws.Range("L136").Value = 10200
$ws.Range("N136").Value = -15300
